# Move a trailing "System" token to the front of the comma-separated
# "Recorded By" list in column G, for every data row on the active sheet.
#
# Example: "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
#          "system, backup@backdoor.com, System" -> "System, system, backup@backdoor.com"
#
# Rows whose value does not end with ", System" (e.g. just "System", or a
# list that doesn't contain "System" at all) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value

    if ($text.EndsWith(", System")) {
        $parts = $text.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        # Last token is "System" (verified by EndsWith check above)
        $rest = $parts[0..($parts.Length - 2)]
        $newParts = @("System") + $rest
        $newText = [string]::Join(", ", $newParts)
        $cell.Value = $newText
    }
}
